$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64 - this shifts the old rows 64-70 down to 65-71,
# matching the target diff (old row64 -> new row65, ..., old row70 -> new row71).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new weekly price record.
$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value = "Los Lagos"
$ws.Cells.Item(64, 4).Value = 45212
$ws.Cells.Item(64, 5).Value = 10
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100101
$ws.Cells.Item(64, 8).Value = "Berries"
$ws.Cells.Item(64, 9).Value = 100101001
$ws.Cells.Item(64, 10).Value = "Arándano (blue)"
$ws.Cells.Item(64, 11).Value = "Sin especificar"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 200
$ws.Cells.Item(64, 14).Value = 13000
$ws.Cells.Item(64, 15).Value = 13000
$ws.Cells.Item(64, 16).Value = 13000
$ws.Cells.Item(64, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(64, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(64, 19).Value = 6500
$ws.Cells.Item(64, 20).Value = 2
